$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 500
$ws.Range("I6").Value = 500
$ws.Range("K6").Value = 1500
$ws.Range("M6").Value = -1388
$ws.Range("H103").Value = 1389.625
$ws.Range("I103").Value = 1929.3334
$ws.Range("J103").Value = 1065.8
$ws.Range("K103").Value = 5788.0002
$ws.Range("L103").Value = 3197.4
$ws.Range("M103").Value = -5202.0002
$ws.Range("N103").Value = -4369.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 964.78314
$ws.Range("I45").Value = 913.4875
$ws.Range("J45").Value = 2332.6667
$ws.Range("K45").Value = 913.4875
$ws.Range("L45").Value = 2332.6667
$ws.Range("M45").Value = -536.4875
$ws.Range("N45").Value = -3086.6667
$ws.Range("H61").Value = 2912.3572
$ws.Range("I61").Value = 2597.923
$ws.Range("J61").Value = 7000
$ws.Range("K61").Value = 2597.923
$ws.Range("L61").Value = 7000
$ws.Range("M61").Value = -2385.923
$ws.Range("N61").Value = -7424
$ws.Range("H74").Value = 1744
$ws.Range("I74").Value = 1931.75
$ws.Range("K74").Value = 1931.75
$ws.Range("M74").Value = -1057.75
$ws.Range("H77").Value = 1744
$ws.Range("I77").Value = 1931.75
$ws.Range("K77").Value = 9658.75
$ws.Range("M77").Value = -5290.75
$ws.Range("H132").Value = 4178.278
$ws.Range("I132").Value = 3093.9092
$ws.Range("J132").Value = 5882.2856
$ws.Range("K132").Value = 9281.7276
$ws.Range("L132").Value = 17646.8568
$ws.Range("M132").Value = -6751.7276
$ws.Range("N132").Value = -22706.8568
$ws.Range("H136").Value = 2912.3572
$ws.Range("I136").Value = 2597.923
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 7793.768999999999
$ws.Range("L136").Value = 21000
$ws.Range("M136").Value = -5243.768999999999
$ws.Range("N136").Value = -26100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2569.6155
$ws.Range("I20").Value = 2450.4167
$ws.Range("K20").Value = 2450.4167
$ws.Range("M20").Value = -2203.4167
$ws.Range("H94").Value = 475.3846
$ws.Range("I94").Value = 475.3846
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 475.3846
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -24.38459999999998
$ws.Range("N94").ClearContents()
$ws.Range("H99").Value = 1667.9166
$ws.Range("I99").Value = 1619.5454
$ws.Range("J99").Value = 2200
$ws.Range("K99").Value = 1619.5454
$ws.Range("L99").Value = 2200
$ws.Range("M99").Value = -121.5454
$ws.Range("N99").Value = -5196
$ws.Range("H134").Value = 3194.8696
$ws.Range("I134").Value = 3084.647
$ws.Range("J134").Value = 3507.1667
$ws.Range("K134").Value = 9253.940999999999
$ws.Range("L134").Value = 10521.5001
$ws.Range("M134").Value = -6718.940999999999
$ws.Range("N134").Value = -15591.5001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1938.0834
$ws.Range("I22").Value = 1627.1666
$ws.Range("J22").Value = 2249
$ws.Range("K22").Value = 1627.1666
$ws.Range("L22").Value = 2249
$ws.Range("M22").Value = -1277.1666
$ws.Range("N22").Value = -2949
$ws.Range("H58").Value = 3892.5715
$ws.Range("I58").Value = 2726.875
$ws.Range("K58").Value = 2726.875
$ws.Range("M58").Value = -2523.875
$ws.Range("H99").Value = 24460.5
$ws.Range("I99").Value = 26086.5
$ws.Range("K99").Value = 26086.5
$ws.Range("M99").Value = -24588.5
$ws.Range("H122").Value = 4633.1665
$ws.Range("I122").Value = 4633.1665
$ws.Range("K122").Value = 13899.4995
$ws.Range("M122").Value = -11449.4995
$ws.Range("H126").Value = 24460.5
$ws.Range("I126").Value = 26086.5
$ws.Range("K126").Value = 78259.5
$ws.Range("M126").Value = -75789.5
$ws.Range("H132").Value = 3761.625
$ws.Range("I132").Value = 3017.182
$ws.Range("J132").Value = 5399.4
$ws.Range("K132").Value = 9051.545999999998
$ws.Range("L132").Value = 16198.2
$ws.Range("M132").Value = -6521.545999999998
$ws.Range("N132").Value = -21258.2
$ws.Range("H134").Value = 2109.8
$ws.Range("I134").Value = 2109.8
$ws.Range("K134").Value = 6329.400000000001
$ws.Range("M134").Value = -3794.400000000001
$ws.Range("H136").Value = 3892.5715
$ws.Range("I136").Value = 2726.875
$ws.Range("K136").Value = 8180.625
$ws.Range("M136").Value = -5630.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7631.077
$ws.Range("I80").Value = 3366
$ws.Range("J80").Value = 8910.6
$ws.Range("K80").Value = 3366
$ws.Range("L80").Value = 8910.6
$ws.Range("M80").Value = -2368
$ws.Range("N80").Value = -10906.6
$ws.Range("H83").Value = 7631.077
$ws.Range("I83").Value = 3366
$ws.Range("J83").Value = 8910.6
$ws.Range("K83").Value = 16830
$ws.Range("L83").Value = 44553
$ws.Range("M83").Value = -11838
$ws.Range("N83").Value = -54537
$ws.Range("H113").Value = 3195.9167
$ws.Range("I113").Value = 3159.182
$ws.Range("K113").Value = 3159.182
$ws.Range("M113").Value = -989.1819999999998
$ws.Range("H122").Value = 1881.9445
$ws.Range("I122").Value = 1905.5834
$ws.Range("J122").Value = 1834.6666
$ws.Range("K122").Value = 5716.7502
$ws.Range("L122").Value = 5503.9998
$ws.Range("M122").Value = -3266.7502
$ws.Range("N122").Value = -10403.9998
$ws.Range("H126").Value = 5166.05
$ws.Range("I126").Value = 4376.8
$ws.Range("J126").Value = 5955.3
$ws.Range("K126").Value = 13130.4
$ws.Range("L126").Value = 17865.9
$ws.Range("M126").Value = -10660.4
$ws.Range("N126").Value = -22805.9
$ws.Range("H132").Value = 3422.4333
$ws.Range("I132").Value = 3158.8635
$ws.Range("J132").Value = 4147.25
$ws.Range("K132").Value = 9476.5905
$ws.Range("L132").Value = 12441.75
$ws.Range("M132").Value = -6946.5905
$ws.Range("N132").Value = -17501.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4028.08
$ws.Range("I7").Value = 3160.7693
$ws.Range("J7").Value = 4967.6665
$ws.Range("K7").Value = 3160.7693
$ws.Range("L7").Value = 4967.6665
$ws.Range("M7").Value = -3048.7693
$ws.Range("N7").Value = -5191.6665
$ws.Range("H126").Value = 4028.08
$ws.Range("I126").Value = 3160.7693
$ws.Range("J126").Value = 4967.6665
$ws.Range("K126").Value = 9482.3079
$ws.Range("L126").Value = 14902.9995
$ws.Range("M126").Value = -7012.3079
$ws.Range("N126").Value = -19842.9995
$ws.Range("H132").Value = 3760.7778
$ws.Range("I132").Value = 3141.4167
$ws.Range("J132").Value = 4999.5
$ws.Range("K132").Value = 9424.250100000001
$ws.Range("L132").Value = 14998.5
$ws.Range("M132").Value = -6894.250100000001
$ws.Range("N132").Value = -20058.5
$ws.Range("H136").Value = 5591.375
$ws.Range("I136").Value = 4704.0625
$ws.Range("K136").Value = 14112.1875
$ws.Range("M136").Value = -11562.1875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 11606.134
$ws.Range("J132").Value = 12572
$ws.Range("L132").Value = 37716
$ws.Range("N132").Value = -42776
$ws.Range("H136").Value = 4086
$ws.Range("I136").Value = 3560.2
$ws.Range("K136").Value = 10680.6
$ws.Range("M136").Value = -8130.599999999999

Write-Output "done"